$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Indiana Pacers vs Milwaukee Bucks"
$ws.Range("B2").Value = "26-04-2024"
$ws.Range("C2").Value = "Indianápolis"

$ws.Range("A3").Value = "Dallas Mavericks vs LA Clippers"
$ws.Range("B3").Value = "26-04-2024"
$ws.Range("C3").Value = "Dallas"

$ws.Range("A4").Value = "Phoenix Suns vs Minnesota Timberwolves"
$ws.Range("B4").Value = "26-04-2024"
$ws.Range("C4").Value = "Phoenix"

$ws.Range("A5:C7").Delete()
